$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text so values like "1.00" / "17.00" are not
# silently coerced into numbers (which would drop the trailing zeros / dots).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '56.349.63'
$ws.Range('E2').Value = '  -1.44%  '

$ws.Range('D3').Value = '2.323.22'
$ws.Range('E3').Value = '  -0.79%  '

$ws.Range('E4').Value = '  +0.18%  '

$ws.Range('D5').Value = '510.97'
$ws.Range('E5').Value = '  -1.79%  '

$ws.Range('D6').Value = '131.44'
$ws.Range('E6').Value = '  -2.72%  '

$ws.Range('E7').Value = '  +0.17%  '

$ws.Range('E8').Value = '  -1.36%  '

$ws.Range('E9').Value = '  -3.21%  '

$ws.Range('E10').Value = '  -0.40%  '

$ws.Range('D11').Value = '5.24'
$ws.Range('E11').Value = '  +0.17%  '

$ws.Range('D12').Value = '0.336'
$ws.Range('E12').Value = '  -1.62%  '

$ws.Range('D13').Value = '2.738.65'
$ws.Range('E13').Value = '  -0.18%  '

$ws.Range('D14').Value = '23.44'
$ws.Range('E14').Value = '  -1.07%  '

$ws.Range('D15').Value = '56.369.71'
$ws.Range('E15').Value = '  -1.14%  '

$ws.Range('E16').Value = '  -2.00%  '

$ws.Range('D17').Value = '2.321.37'
$ws.Range('E17').Value = '  -0.97%  '

$ws.Range('D18').Value = '10.42'
$ws.Range('E18').Value = '  -0.66%  '

$ws.Range('D19').Value = '325.09'
$ws.Range('E19').Value = '  +0.26%  '

$ws.Range('D20').Value = '4.11'
$ws.Range('E20').Value = '  -2.79%  '

$ws.Range('E21').Value = '  +2.89%  '

$ws.Range('E22').Value = '  +0.00%  '

$ws.Range('D23').Value = '61.56'
$ws.Range('E23').Value = '  +0.78%  '

$ws.Range('D24').Value = '8.74'
$ws.Range('E24').Value = '  +11.14%  '

$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.62%  '

$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').Value = '  -0.79%  '

$ws.Range('D27').Value = '1.30'
$ws.Range('E27').Value = '  +1.94%  '

$ws.Range('D28').Value = '167.33'
$ws.Range('E28').Value = '  -1.98%  '

$ws.Range('D29').Value = '1.67'
$ws.Range('E29').Value = '  -2.67%  '

$ws.Range('E30').Value = '  -4.41%  '

$ws.Range('E31').Value = '  -1.54%  '

$ws.Range('E32').Value = '  -0.08%  '

$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('E34').Value = '  +0.07%  '

$ws.Range('E35').Value = '  +0.37%  '

$ws.Range('E36').Value = '  -4.06%  '

$ws.Range('E37').Value = '  -3.19%  '

$ws.Range('D38').Value = '38.40'
$ws.Range('E38').Value = '  +1.27%  '

$ws.Range('D39').Value = '1.55'
$ws.Range('E39').Value = '  -0.52%  '

$ws.Range('D40').Value = '148.92'
$ws.Range('E40').Value = '  +8.69%  '

$ws.Range('E41').Value = '  -1.74%  '

$ws.Range('E42').Value = '  -1.40%  '

$ws.Range('D43').Value = '276.19'
$ws.Range('E43').Value = '  -0.87%  '

$ws.Range('E44').Value = '  -2.87%  '

$ws.Range('D45').Value = '0.0925'
$ws.Range('E45').Value = '  -1.02%  '

$ws.Range('D46').Value = '0.0493'
$ws.Range('E46').Value = '  -2.28%  '

$ws.Range('E47').Value = '  -1.41%  '

$ws.Range('E48').Value = '  +2.43%  '

$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D49').Value = '0.377'
$ws.Range('E49').Value = '  -0.68%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').Value = '0.0213'
$ws.Range('E50').Value = '  -2.56%  '

$ws.Range('D51').Value = '17.00'
$ws.Range('E51').Value = '  +0.95%  '

# Restore the default (unstyled) look for the Price column now that the text
# values are safely stored as strings.
$ws.Range("D2:D51").Style = "Normal"
